# Insert a new data row before the existing row 690 ("Hortaliza, Vega Modelo de
# Temuco - Sandia" sheet). Excel will shift rows 690-738 down to 691-739 and
# extend the used range (dimension) to A1:R739 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(690).Insert()

$row = 690
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 45106
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100112028
$ws.Cells.Item($row, 7).Value = "Sandia"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 150
$ws.Cells.Item($row, 11).Value = 3600
$ws.Cells.Item($row, 12).Value = 3600
$ws.Cells.Item($row, 13).Value = 3600
$ws.Cells.Item($row, 14).Value = "`$/unidad"
$ws.Cells.Item($row, 15).Value = "Brasil"
$ws.Cells.Item($row, 16).Value = 3600
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
